$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 43 (pushes existing rows 43:100 down to 45:102)
$ws.Rows("43:44").Insert()

# New week's price report for "Pepino ensalada" - Comercializadora del Agro de Limarí
# Row 43: Calidad "Primera"
$ws.Range("A43").Value = 2
$ws.Range("B43").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C43").Value = "Coquimbo"
$ws.Range("D43").Value = 44546
$ws.Range("E43").Value = 4
$ws.Range("F43").Value = 100112043
$ws.Range("G43").Value = "Pepino ensalada"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 400
$ws.Range("K43").Value = 6000
$ws.Range("L43").Value = 7000
$ws.Range("M43").Value = 6500
$ws.Range("N43").Value = "`$/caja 70 unidades"
$ws.Range("O43").Value = "Provincia de Limarí"
$ws.Range("P43").Value = 93
$ws.Range("Q43").Value = 70
$ws.Range("R43").Value = "Hortaliza"

# Row 44: Calidad "Segunda"
$ws.Range("A44").Value = 2
$ws.Range("B44").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C44").Value = "Coquimbo"
$ws.Range("D44").Value = 44546
$ws.Range("E44").Value = 4
$ws.Range("F44").Value = 100112043
$ws.Range("G44").Value = "Pepino ensalada"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Segunda"
$ws.Range("J44").Value = 240
$ws.Range("K44").Value = 4000
$ws.Range("L44").Value = 5000
$ws.Range("M44").Value = 4500
$ws.Range("N44").Value = "`$/caja 100 unidades"
$ws.Range("O44").Value = "Provincia de Limarí"
$ws.Range("P44").Value = 45
$ws.Range("Q44").Value = 100
$ws.Range("R44").Value = "Hortaliza"
